$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# cryptos list, as pulled by the scheduled GitHub Actions job.
# Cells whose new text is unambiguously numeric are forced to keep their
# original text representation (e.g. "1.000", "0.7063") by setting the
# cell number format to Text before assigning the value, matching the
# inline-string storage used throughout this sheet.
$ws.Range("D2").Value = '29.215.31'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.855.94'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7063'
$ws.Range("E5").Value = '  +2.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.03'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08014'
$ws.Range("E8").Value = '  +3.76%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08176'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.182'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Value = '1.810.49'
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7017'
$ws.Range("E14").Value = '  -3.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.52'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '29.111.68'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.793'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007879'
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.18'
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9988'
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '2.041.47'
$ws.Range("E23").Value = '  -3.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.454'
$ws.Range("E24").Value = '  -1.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.76'
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.896'
$ws.Range("E26").Value = '  -0.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1432'
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.09'
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.919'
$ws.Range("E29").Value = '  -2.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.418'
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("E32").Value = '  -2.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.021'
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05186'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7122'
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9977'
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.639'
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01849'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.722'
$ws.Range("E40").Value = '  +1.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9309'
$ws.Range("E41").Value = '  +1.70%  '
$ws.Range("D42").Value = '1.132.37'
$ws.Range("E42").Value = '  +4.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.908'
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.95'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9999'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.44'
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5328'
$ws.Range("E48").Value = '  -4.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.757'
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.160'
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.938'
$ws.Range("E51").Value = '  -0.62%  '
